$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.041, 0.149, -0.04, -49, 5, -25),
    @(-0.426, -0.471, 0.553, -93, -10, 33),
    @(0.08799999999999999, -0.613, 0.133, -78, -45, -105),
    @(-0.038, 0.16, 0.174, 26, -12, -23),
    @(0.113, -0.127, -0.016, -63, -1, -154),
    @(0.424, -0.042, -0.271, -141, 57, -109),
    @(0.336, -0.375, 0.09, 16, 7, -108),
    @(0.129, -0.434, 0.5, 150, -17, -63),
    @(-0.217, -0.165, 0.09, -128, 57, -177),
    @(0.135, -0.136, 0.451, -73, -16, -104)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}
